# Actualización automática 2025-10-08 17:30:10
$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
# New octubre PORCELANATO sale registered for F.V - AREA ANDINA S.A. (row 12)
$wsGrupo.Range("M12").Value = 626.6900000000001
# Count of clients with PORCELANATO sales goes from 2 to 3 (out of 24)
$wsGrupo.Range("M26").Value = "3 de 24"

# --- Sheet "VENTA MENSUAL" ---
# October (octubre) sale for the same client
$wsMensual.Range("F12").Value = 626.6900000000001
# Updated October total
$wsMensual.Range("F26").Value = 5756.16

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
# PORCELANATO row (row 12): VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumplimiento.Range("D12").Value = 7148.19
$wsCumplimiento.Range("E12").Value = 20806.79
$wsCumplimiento.Range("F12").Value = 0.255703634915854

# TOTAL row (row 14): VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumplimiento.Range("D14").Value = 5756.16
$wsCumplimiento.Range("E14").Value = 36447.22110009468
$wsCumplimiento.Range("F14").Value = 0.1363909679735846
